# Add the "2022-Q3" quarterly report sheet, positioned right after "总计",
# and record its summary figures on the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the "总计" (totals) summary sheet: insert a new row for
#    2022-Q3 right under the header, pushing the existing quarters down.
# ---------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows(2).Insert()

# Insert() clones formatting from the row below into the blank row;
# strip that before writing so the new data cells end up unstyled,
# matching the rest of the sheet's body rows.
$totals.Range("B2:D2").ClearFormats()

# Column A carries the running index style, so copy that look from
# the row just beneath (which used to be row 2).
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 7
$totals.Range("D2").Value = 0.51

# Renumber the running index in column A for the quarters that got
# pushed down a row.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3

# ---------------------------------------------------------------
# 2. Create the new "2022-Q3" sheet right after "总计" and fill it
#    with the quarter's fund holdings.
# ---------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totals)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Borders.LineStyle = 1

# Fund code / numeric-text columns must stay text so leading zeros and
# trailing zeros in the source data are preserved exactly.
$q3.Range("B2:G8").NumberFormat = "@"

$funds = @(
    @("010861", "长信企业优选一年持有期灵活配置混合", "8.09", "81.28", "3.90", "0.3155", 2),
    @("005589", "长信企业精选两年定期开放灵活配置混合", "2.06", "81.65", "3.85", "0.0793", 1),
    @("014356", "长信企业成长三年持有混合A", "1.70", "82.49", "3.88", "0.0660", 1),
    @("014357", "长信企业成长三年持有混合C", "0.87", "82.49", "3.88", "0.0338", 1),
    @("007294", "长信利信灵活配置混合E", "0.49", "57.23", "3.01", "0.0147", 4),
    @("007293", "长信利信灵活配置混合C", "0.04", "57.23", "3.01", "0.0012", 4),
    @("519949", "长信利信灵活配置混合A", "0.01", "57.23", "3.01", "0.0003", 4)
)

$r = 2
foreach ($fund in $funds) {
    $q3.Range("A$r").Value = $r - 2
    $q3.Range("B$r").Value = $fund[0]
    $q3.Range("C$r").Value = $fund[1]
    $q3.Range("D$r").Value = $fund[2]
    $q3.Range("E$r").Value = $fund[3]
    $q3.Range("F$r").Value = $fund[4]
    $q3.Range("G$r").Value = $fund[5]
    $q3.Range("H$r").Value = $fund[6]
    $r = $r + 1
}

$q3.Range("A2:A8").Font.Bold = $true
$q3.Range("A2:A8").HorizontalAlignment = -4108
$q3.Range("A2:A8").VerticalAlignment = -4160
$q3.Range("A2:A8").Borders.LineStyle = 1

$q3.Range("A1").Select()
